$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "countryLang" column (F) with header and per-row values
$ws.Range("F1").Value = "countryLang"

$ws.Range("F2").Value = "te-IN"
$ws.Range("F3").Value = "en-IN"
$ws.Range("F4").Value = "en-KE"
$ws.Range("F5").Value = "kn-IN"

# Match header style (bold) used by the rest of row 1
$ws.Range("F1").Font.Bold = $true

# Update selection to the last edited cell, matching the saved view state
$ws.Range("F5").Select()
